# Auto-generated update: add data for 2025-09-14 to violent-crime-full-year workbook
# Updates the "2025" (column L) year-to-date counts across the Citywide Totals, By
# Neighborhood summary, and every individual neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4739  # was 4714
$ws.Range("L3").Value = 5115  # was 5082
$ws.Range("L4").Value = 1258  # was 1254
$ws.Range("L5").Value = 301  # was 299
$ws.Range("L6").Value = 4338  # was 4316
$ws.Range("L7").Value = 15751  # was 15665

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 120  # was 119
$ws.Range("L7").Value = 516  # was 513
$ws.Range("L8").Value = 1044  # was 1040
$ws.Range("L10").Value = 105  # was 104
$ws.Range("L11").Value = 255  # was 254
$ws.Range("L15").Value = 115  # was 114
$ws.Range("L18").Value = 111  # was 110
$ws.Range("L19").Value = 433  # was 430
$ws.Range("L23").Value = 172  # was 169
$ws.Range("L27").Value = 140  # was 139
$ws.Range("L29").Value = 864  # was 863
$ws.Range("L30").Value = 75  # was 74
$ws.Range("L31").Value = 157  # was 156
$ws.Range("L33").Value = 720  # was 713
$ws.Range("L36").Value = 204  # was 202
$ws.Range("L37").Value = 582  # was 576
$ws.Range("L40").Value = 42  # was 41
$ws.Range("L41").Value = 71  # was 70
$ws.Range("L42").Value = 515  # was 512
$ws.Range("L43").Value = 116  # was 114
$ws.Range("L44").Value = 111  # was 110
$ws.Range("L47").Value = 109  # was 108
$ws.Range("L51").Value = 194  # was 193
$ws.Range("L52").Value = 315  # was 314
$ws.Range("L54").Value = 330  # was 329
$ws.Range("L57").Value = 56  # was 55
$ws.Range("L60").Value = 100  # was 99
$ws.Range("L63").Value = 44  # was 41
$ws.Range("L65").Value = 305  # was 304
$ws.Range("L66").Value = 39  # was 38
$ws.Range("L67").Value = 541  # was 539
$ws.Range("L73").Value = 123  # was 122
$ws.Range("L75").Value = 58  # was 57
$ws.Range("L78").Value = 210  # was 209
$ws.Range("L79").Value = 414  # was 412
$ws.Range("L83").Value = 349  # was 345
$ws.Range("L84").Value = 153  # was 151
$ws.Range("L85").Value = 811  # was 804
$ws.Range("L90").Value = 159  # was 158
$ws.Range("L91").Value = 214  # was 211
$ws.Range("L93").Value = 82  # was 81
$ws.Range("L99").Value = 273  # was 270
$ws.Range("L100").Value = 27  # was 25
$ws.Range("L101").Value = 15751  # was 15665

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 175  # was 174
$ws.Range("L3").Value = 171  # was 169
$ws.Range("L7").Value = 516  # was 513

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 78  # was 77
$ws.Range("L7").Value = 255  # was 254

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 245  # was 241
$ws.Range("L3").Value = 328  # was 326
$ws.Range("L6").Value = 169  # was 168
$ws.Range("L7").Value = 811  # was 804

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 84  # was 83
$ws.Range("L7").Value = 315  # was 314

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 302  # was 300
$ws.Range("L6").Value = 276  # was 274
$ws.Range("L7").Value = 1044  # was 1040

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 107  # was 105
$ws.Range("L3").Value = 139  # was 138
$ws.Range("L6").Value = 81  # was 80
$ws.Range("L7").Value = 349  # was 345

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 195  # was 194
$ws.Range("L3").Value = 247  # was 242
$ws.Range("L6").Value = 220  # was 219
$ws.Range("L7").Value = 720  # was 713

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 171  # was 170
$ws.Range("L3").Value = 202  # was 201
$ws.Range("L4").Value = 32  # was 31
$ws.Range("L6").Value = 160  # was 157
$ws.Range("L7").Value = 582  # was 576

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L6").Value = 82  # was 81
$ws.Range("L7").Value = 305  # was 304

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 114  # was 111
$ws.Range("L7").Value = 273  # was 270

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L3").Value = 18  # was 17
$ws.Range("L7").Value = 75  # was 74

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L3").Value = 41  # was 40
$ws.Range("L7").Value = 157  # was 156

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L4").Value = 39  # was 38
$ws.Range("L6").Value = 124  # was 123
$ws.Range("L7").Value = 541  # was 539

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 49  # was 47
$ws.Range("L7").Value = 153  # was 151

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 80  # was 79
$ws.Range("L7").Value = 330  # was 329

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 258  # was 257
$ws.Range("L3").Value = 327  # was 326
$ws.Range("L5").Value = 14  # was 15
$ws.Range("L7").Value = 864  # was 863

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 152  # was 150
$ws.Range("L3").Value = 135  # was 134
$ws.Range("L7").Value = 433  # was 430

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L6").Value = 28  # was 27
$ws.Range("L7").Value = 111  # was 110

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L6").Value = 23  # was 22
$ws.Range("L7").Value = 120  # was 119

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L2").Value = 23  # was 22
$ws.Range("L7").Value = 71  # was 70

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 172  # was 171
$ws.Range("L4").Value = 38  # was 37
$ws.Range("L5").Value = 14  # was 13
$ws.Range("L7").Value = 515  # was 512

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 45  # was 44
$ws.Range("L7").Value = 105  # was 104

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 64  # was 63
$ws.Range("L7").Value = 210  # was 209

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 65  # was 63
$ws.Range("L6").Value = 47  # was 46
$ws.Range("L7").Value = 172  # was 169

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 76  # was 74
$ws.Range("L6").Value = 26  # was 25
$ws.Range("L7").Value = 214  # was 211

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L5").Value = 13  # was 12
$ws.Range("L6").Value = 88  # was 87
$ws.Range("L7").Value = 414  # was 412

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 40  # was 39
$ws.Range("L7").Value = 111  # was 110

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 76  # was 75
$ws.Range("L3").Value = 60  # was 59
$ws.Range("L7").Value = 204  # was 202

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L6").Value = 26  # was 25
$ws.Range("L7").Value = 82  # was 81

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 15  # was 13
$ws.Range("L7").Value = 27  # was 25

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 36  # was 35
$ws.Range("L7").Value = 109  # was 108

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 41  # was 40
$ws.Range("L7").Value = 115  # was 114

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L2").Value = 11  # was 10
$ws.Range("L7").Value = 39  # was 38

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 38  # was 37
$ws.Range("L7").Value = 123  # was 122

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L2").Value = 39  # was 38
$ws.Range("L7").Value = 140  # was 139

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L3").Value = 21  # was 20
$ws.Range("L7").Value = 58  # was 57

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 47  # was 46
$ws.Range("L7").Value = 159  # was 158

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 64  # was 63
$ws.Range("L7").Value = 194  # was 193

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L3").Value = 16  # was 15
$ws.Range("L7").Value = 56  # was 55

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L5").Value = 5  # was 4
$ws.Range("L7").Value = 100  # was 99

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 40  # was 38
$ws.Range("L7").Value = 116  # was 114

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 19  # was 18
$ws.Range("L7").Value = 42  # was 41
